# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described in the commit diff (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force Text storage so values such as "1.00" / "0.999"
# keep their exact digits instead of being normalised by Excel's automatic
# number detection. NumberFormat is restored to the default "Normal" style
# right after the write so no new cell style is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.461.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.439.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.419.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.829"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.787.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.16"
$ws.Range("D49").Style = "Normal"

# --- 1-hour volume/change column (E): values keep their leading/trailing
# double-space padding, which also keeps Excel from reinterpreting the
# "+x.xx%" / "-x.xx%" text as a percentage number.

$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("E18").Value = "  -4.88%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -7.15%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("E51").Value = "  -3.21%  "
